# Apply cryptos list price/volume updates (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.353.22"
$ws.Range("E2").Value = "  -0.68%  "

$ws.Range("D3").Value = "1.797.40"
$ws.Range("E3").Value = "  -1.12%  "

$ws.Range("E4").Value = "  -0.33%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "225.14"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -1.51%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "0.599"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +3.60%  "

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "1.00"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  -0.27%  "

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "36.28"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  +3.50%  "

$ws.Range("E9").Value = "  -3.29%  "

$ws.Range("E10").Value = "  -3.11%  "

$ws.Range("E11").Value = "  +1.14%  "

$ws.Range("D12").Value = "2.056.97"
$ws.Range("E12").Value = "  -1.15%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "11.23"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  -1.23%  "

$ws.Range("D14").Value = "1.802.34"
$ws.Range("E14").Value = "  -1.16%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "0.630"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  -2.81%  "

$ws.Range("D16").Value = "34.309.06"
$ws.Range("E16").Value = "  -0.82%  "

$ws.Range("E17").Value = "  +1.16%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "68.62"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  -1.20%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "245.78"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "0.0₃0772"
$ws.Range("E20").Value = "  -3.70%  "

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "11.31"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  -2.22%  "

$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("E23").Value = "  -2.88%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "2.21"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +4.83%  "

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "170.40"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  -0.94%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "7.86"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  +4.09%  "

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "17.32"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  +2.97%  "

$ws.Range("E28").Value = "  +1.40%  "

$ws.Range("E29").Value = "  -0.33%  "

$ws.Range("E30").Value = "  -1.78%  "

$ws.Range("E31").Value = "  -1.61%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "3.89"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -2.38%  "

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0513"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  -3.32%  "

$ws.Range("E34").Value = "  -4.39%  "

$ws.Range("D35").Value = "1.360.84"
$ws.Range("E35").Value = "  -3.02%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "0.646"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  -5.05%  "

$ws.Range("E37").Value = "  -1.70%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "2.35"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -8.06%  "

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0186"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  -2.77%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "2.42"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  +0.47%  "

$ws.Range("E41").Value = "  -2.64%  "

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "80.72"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  -3.00%  "

$ws.Range("E43").Value = "  -1.94%  "

$ws.Range("E44").Value = "  +4.91%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "13.19"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -4.84%  "

$ws.Range("E46").Value = "  -3.50%  "

$ws.Range("D47").Value = "1.958.09"
$ws.Range("E47").Value = "  -1.17%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "5.74"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  -5.16%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "0.999"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  -0.28%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "101.73"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  -3.63%  "

$ws.Range("D51").Value = "0.0₆0120"
$ws.Range("E51").Value = "  -8.24%  "
